# Applies the 2025-12-25 10:16:55 sync of attendance_reports data:
#  1. Updates the "Class Statistics" summary metrics (Recorded/Missing Sessions,
#     Coverage %, Average Attendance %) to reflect the newly recorded sessions.
#  2. Normalizes the "Recorded By" ordering from "dnasr281@gmail.com, System"
#     to "System, dnasr281@gmail.com" everywhere it occurs.
#  3. Updates the "Group Statistics" rollup (Recorded/Missing counts and the
#     Coverage %/Avg Attendance % columns) for the B1A1 group's later sessions.
#  4. Marks the six 25/12/2025 (B1D1/B1D2/B1E1/B1E2/B1F1/B1F2) sessions as
#     Recorded now that attendance has been taken, copying the "Recorded" row
#     formatting (green fill) over the previous "Not Recorded" (pink) formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Class Statistics summary block (column L)
#    Note: percentage-looking values are written with a leading apostrophe
#    so Excel stores them as literal text ("53.8%") rather than silently
#    re-interpreting them as a numeric percentage (0.538 formatted as %),
#    matching the original report's text-based percentage columns.
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 171          # Recorded Sessions
$ws.Range("L7").Value = 3            # Missing Sessions
$ws.Range("L9").Value = "'53.8%"     # Coverage %
$ws.Range("L10").Value = "'75.2%"    # Average Attendance %

# ---------------------------------------------------------------------------
# 2) "Recorded By" ordering swap in column G (applies to every row where it
#    currently reads "dnasr281@gmail.com, System")
# ---------------------------------------------------------------------------
$gSwapRows = @(8,9,10,34,35,36,60,61,62,86,87,88,112,113,114,138,139,140,
               164,167,170,191,194,197,218,221,224,245,248,251,272,275,278,
               299,302,305)

foreach ($r in $gSwapRows) {
    $cell = $ws.Range("G$r")
    if ($cell.Value() -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
}

# ---------------------------------------------------------------------------
# 3) Group Statistics rollup for B1A1 (rows 21-26), columns O/P/R/S
# ---------------------------------------------------------------------------
$groupStatsChanges = @{
    21 = @{ O = 15; P = 0; R = "'55.6%"; S = "'78.0%" }
    22 = @{ O = 15; P = 0; R = "'55.6%"; S = "'76.7%" }
    23 = @{ O = 15; P = 0; R = "'55.6%"; S = "'80.0%" }
    24 = @{ O = 14; P = 1; R = "'51.9%"; S = "'71.2%" }
    25 = @{ O = 15; P = 0; R = "'55.6%"; S = "'70.3%" }
    26 = @{ O = 15; P = 0; R = "'55.6%"; S = "'62.3%" }
}

foreach ($r in $groupStatsChanges.Keys) {
    $vals = $groupStatsChanges[$r]
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("R$r").Value = $vals.R
    $ws.Range("S$r").Value = $vals.S
}

# ---------------------------------------------------------------------------
# 4) Flip six "Not Recorded" sessions (25/12/2025) to "Recorded" - copy the
#    green "Recorded" formatting from a neighboring recorded row, then fill
#    in the Recorded By / Students / Status values.
# ---------------------------------------------------------------------------
$rowRewrites = @{
    172 = @{ H = "19/23"; SourceRow = 170 }
    199 = @{ H = "26/30"; SourceRow = 197 }
    226 = @{ H = "18/25"; SourceRow = 224 }
    253 = @{ H = "21/28"; SourceRow = 251 }
    280 = @{ H = "22/26"; SourceRow = 278 }
    307 = @{ H = "16/29"; SourceRow = 305 }
}

foreach ($r in $rowRewrites.Keys) {
    $info = $rowRewrites[$r]
    $src = $info.SourceRow

    $ws.Range("A$($src):I$($src)").Copy()
    $ws.Range("A$($r):I$($r)").PasteSpecial(-4122)   # xlPasteFormats

    $ws.Range("G$r").Value = "dnasr281@gmail.com"
    $ws.Range("H$r").Value = $info.H
    $ws.Range("I$r").Value = "Recorded"
}

$excel.CutCopyMode = 0
